# Apply updated cryptocurrency price/volume data (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as plain text (locale-formatted, e.g. "37.511.14" or
# trailing zeros like "15.00"). Mark the cells whose new value would otherwise
# be auto-sniffed as a number ("@" = Text) BEFORE writing them, so Excel keeps
# the exact text representation instead of collapsing it to a float.
$textRanges = @(
    "D5:D7",
    "D10:D11",
    "D13:D14",
    "D16:D17",
    "D20",
    "D22:D24",
    "D27:D30",
    "D32",
    "D34",
    "D37:D40",
    "D43:D44",
    "D46:D48",
    "D50:D51"
)
foreach ($r in $textRanges) {
    $ws.Range($r).NumberFormat = "@"
}

$updates = @(
    @{ Row = 2; D = "37.511.14"; E = "  +5.80%  " },
    @{ Row = 3; D = "2.060.67" },
    @{ Row = 4; E = "  -0.18%  " },
    @{ Row = 5; D = "253.02"; E = "  +3.60%  " },
    @{ Row = 6; D = "0.652"; E = "  +3.01%  " },
    @{ Row = 7; D = "66.55"; E = "  +16.62%  " },
    @{ Row = 8; E = "  -0.10%  " },
    @{ Row = 9; E = "  +6.31%  " },
    @{ Row = 10; D = "59.27"; E = "  +0.12%  " },
    @{ Row = 11; D = "0.0770"; E = "  +5.46%  " },
    @{ Row = 12; E = "  +1.56%  " },
    @{ Row = 13; D = "0.915"; E = "  -2.12%  " },
    @{ Row = 14; D = "15.00"; E = "  +6.11%  " },
    @{ Row = 15; D = "2.360.52"; E = "  +4.33%  " },
    @{ Row = 16; D = "5.59"; E = "  +6.97%  " },
    @{ Row = 17; D = "20.73"; E = "  +19.16%  " },
    @{ Row = 18; D = "2.057.86"; E = "  +3.85%  " },
    @{ Row = 19; D = "37.349.69"; E = "  +5.60%  " },
    @{ Row = 20; D = "74.10"; E = "  +4.86%  " },
    @{ Row = 21; D = "0.0₃0880"; E = "  +4.91%  " },
    @{ Row = 22; D = "5.47"; E = "  +6.56%  " },
    @{ Row = 23; D = "240.63"; E = "  +3.71%  " },
    @{ Row = 24; D = "2.66"; E = "  +5.57%  " },
    @{ Row = 25; E = "  -0.09%  " },
    @{ Row = 26; E = "  +3.28%  " },
    @{ Row = 27; D = "9.68"; E = "  +6.89%  " },
    @{ Row = 28; D = "161.80"; E = "  -0.75%  " },
    @{ Row = 29; D = "20.02"; E = "  +4.10%  " },
    @{ Row = 30; D = "5.30"; E = "  +10.00%  " },
    @{ Row = 31; E = "  +3.38%  " },
    @{ Row = 32; D = "0.112"; E = "  +22.16%  " },
    @{ Row = 33; E = "  +6.51%  " },
    @{ Row = 34; D = "4.79"; E = "  +13.16%  " },
    @{ Row = 35; E = "  +4.82%  " },
    @{ Row = 36; E = "  +5.99%  " },
    @{ Row = 37; D = "6.24"; E = "  +22.54%  " },
    @{ Row = 38; D = "0.999"; E = "  -0.29%  " },
    @{ Row = 39; D = "1.83"; E = "  +4.13%  " },
    @{ Row = 40; D = "3.05"; E = "  +35.55%  " },
    @{ Row = 41; E = "  +16.88%  " },
    @{ Row = 42; E = "  +4.41%  " },
    @{ Row = 43; D = "3.00"; E = "  +5.16%  " },
    @{ Row = 44; D = "1.16"; E = "  +6.92%  " },
    @{ Row = 45; E = "  +4.66%  " },
    @{ Row = 46; D = "17.13"; E = "  +7.64%  " },
    @{ Row = 47; D = "95.35"; E = "  +4.55%  " },
    @{ Row = 48; D = "7.91"; E = "  +5.46%  " },
    @{ Row = 49; D = "1.421.07"; E = "  +3.15%  " },
    @{ Row = 50; D = "2.95"; E = "  +2.40%  " },
    @{ Row = 51; D = "46.89"; E = "  +3.42%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}

Write-Output "Updated $($updates.Count) rows"
